$d = $word.ActiveDocument

$r1 = $d.Content
$r1.Find.Execute("What else would you include?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$delStart = $r1.End

$r3 = $d.Content
$r3.Find.Execute("https://lnkd.in/d4p3xbuv", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$delEnd = $r3.End

Write-Host "delStart=$delStart delEnd=$delEnd"

$delRange = $d.Range($delStart, $delEnd)
Write-Host "About to delete:"
Write-Host $delRange.Text
$delRange.Delete()
Write-Host "Deleted."
